$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header labels: add missing Spanish accents
$ws.Range("D1").Value = "Correo electrónico (*)"
$ws.Range("F1").Value = "Subárea 1 (*)"
$ws.Range("G1").Value = "Subárea 2"
$ws.Range("H1").Value = "Subárea 3"

# Update the view: scroll so column J is the top-left visible column,
# and select Q1 as the active cell
$ws.Range("Q1").Select()
$excel.ActiveWindow.ScrollColumn = 10
